# Regenerate the handback report: refresh "Latest Handback DateTime" (column K)
# for the first data row (4ec2308a-f580-44fc-80a2-02fa039d9313 item) on both the
# "zh-cn" and "de-de" language sheets, reflecting a newer handback timestamp.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-09-07 06:38:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-07 06:38:52"
